$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{Row=126; Date="01-07-2021"; B=11601; C=35270; D=215; E=25235},
    @{Row=127; Date="02-07-2021"; B=12197; C=35258; D=215; E=25114},
    @{Row=128; Date="05-07-2021"; B=13370; C=35341; D=215; E=25096},
    @{Row=129; Date="06-07-2021"; B=13496; C=35261; D=215; E=25308},
    @{Row=130; Date="07-07-2021"; B=14069; C=35248; D=215; E=25319},
    @{Row=131; Date="08-07-2021"; B=14475; C=34767; D=215; E=25366},
    @{Row=132; Date="09-07-2021"; B=15200; C=34705; D=215; E=25233},
    @{Row=133; Date="12-07-2021"; B=14384; C=35216; D=215; E=25447},
    @{Row=134; Date="13-07-2021"; B=14637; C=35182; D=215; E=25585},
    @{Row=135; Date="14-07-2021"; B=14820; C=35120; D=215; E=25480},
    @{Row=136; Date="15-07-2021"; B=13843; C=34690; D=215; E=25630},
    @{Row=137; Date="19-07-2021"; B=14079; C=35361; D=215; E=25833},
    @{Row=138; Date="20-07-2021"; B=14964; C=35370; D=215; E=27818},
    @{Row=139; Date="21-07-2021"; B=15484; C=35408; D=215; E=27798},
    @{Row=140; Date="22-07-2021"; B=15342; C=35384; D=215; E=27840},
    @{Row=141; Date="23-07-2021"; B=14036; C=35488; D=215; E=27899},
    @{Row=142; Date="26-07-2021"; B=13663; C=35496; D=215; E=25678},
    @{Row=143; Date="27-07-2021"; B=13521; C=35536; D=215; E=26019},
    @{Row=144; Date="28-07-2021"; B=13689; C=35580; D=215; E=26056},
    @{Row=145; Date="29-07-2021"; B=13411; C=35528; D=215; E=26089},
    @{Row=146; Date="30-07-2021"; B=12671; C=35239; D=213; E=26787}
)

foreach ($item in $rows) {
    $r = $item.Row
    $dateCell = $ws.Cells.Item($r, 1)
    $dateCell.Formula = '="' + $item.Date + '"'
    $dateCell.Copy()
    $dateCell.PasteSpecial(-4163)

    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
}

$excel.CutCopyMode = 0
